# Rename the first sheet from "1" to "Soal" and make it the active/selected
# sheet (it was previously sheet "3" that was selected). Activating it here
# moves the "tabSelected" flag from sheet "3" to "Soal" and updates the
# workbook's active-tab bookkeeping to point at the first sheet.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Soal"
$ws1.Activate()
